# Sort the raw_data_kff table (Table1) by column C ("Implemented Expansion On"),
# ascending, keeping the header row (row 3) in place. This reproduces the
# "Added 2006-13 crosstab ..." commit's reordering of rows 4:55 driven by a
# sort on the "Implemented Expansion On" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("raw_data_kff")

# The data lives inside an Excel Table (ListObject) anchored at A3:C55.
$lo = $ws.ListObjects.Item("Table1")

$sortKeyRange = $ws.Range("C3:C55")

$lo.Sort.SortFields.Clear()
$lo.Sort.SortFields.Add($sortKeyRange, 0, 1)
$lo.Sort.Header = 1
$lo.Sort.Apply()
